# Insert a new data row at row 31 (pushing existing rows 31..120 down to
# 32..121) and populate it with a new "Arveja Verde" price observation for
# Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31 and below down by one row.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C31").Value = 'Ñuble'
$ws.Range("D31").Value2 = 45251
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = 100112022
$ws.Range("G31").Value = 'Arveja Verde'
$ws.Range("H31").Value = 'Sin especificar'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 50
$ws.Range("K31").Value = 25000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 25000
$ws.Range("N31").Value = '$/saco 25 kilos'
$ws.Range("O31").Value = 'Región de Ñuble'
$ws.Range("P31").Value = 1000
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = 'Hortaliza'
